# Insert a new weekly record for "Feria Lagunitas de Puerto Montt" (Acelga)
# at row 164, pushing the existing rows 164:178 down to 165:179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 164:178 down by one row, preserving all
# formatting/styles of the inserted row (copied from the row above).
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A164").Value = 4
$ws.Range("B164").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C164").Value = "Los Lagos"
$ws.Range("D164").Value = 44714
$ws.Range("E164").Value = 10
$ws.Range("F164").Value = 100112009
$ws.Range("G164").Value = "Acelga"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 30
$ws.Range("K164").Value = 12000
$ws.Range("L164").Value = 12000
$ws.Range("M164").Value = 12000
$ws.Range("N164").Value = "$/docena de atados (12 kilos)"
$ws.Range("O164").Value = "Región de La Araucanía"
$ws.Range("P164").Value = 1000
$ws.Range("Q164").Value = 12
$ws.Range("R164").Value = "Hortaliza"
